$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 44-46: coin name/link/price/volume rotated
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.69'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5314'
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.174'
$ws.Range("E46").Value = '  -4.70%  '

# Remaining price/volume updates
$ws.Range("D2").Value = '28.651.40'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.867.17'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.83'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4626'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3912'
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07904'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9704'
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.30'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").Value = '1.860.68'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.727'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.938'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06955'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.35'
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001005'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.96'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '28.651.69'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.321'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.07'
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.118'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").Value = '2.024.27'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.52'
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.32'
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.749'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.997'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.29'
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9338'
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.323'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.342'
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.352'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05832'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02120'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.152'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.878'
$ws.Range("E39").Value = '  +3.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5657'
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.943'
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1779'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07248'
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("E47").Value = '  -8.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.848'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.27'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.343'
$ws.Range("E51").Value = '  +0.82%  '
